# Apply the "class enrollment" update described by the commit:
#   "everything works but lowercase class names on del"
#
# On the "Timeslot Information" sheet, row 5 (student id 1.4447799E7 / admin10)
# now has:
#   - Monday   (col B) slot 1  -> English  (was null)
#   - Tuesday  (col C) slot 11 -> Java     (was null)
#   - Thursday (col E) slot 10 -> Python   (was null)
#
# On the "Class Information" sheet, the corresponding classes' Current
# Occupancy (col G) counts are updated:
#   - English (row 2)  0 -> 1
#   - Python  (row 11) 0 -> 2
#   - Java    (row 12) 1 -> 2

$wb = $excel.ActiveWorkbook

$timeslots = $wb.Worksheets.Item("Timeslot Information")
$timeslots.Cells.Item(5, 2).Value = "English,null,null,null,Biology,null,null,null,null,null,Java,null"
$timeslots.Cells.Item(5, 3).Value = "null,null,null,null,null,null,null,null,null,null,Java,null"
$timeslots.Cells.Item(5, 5).Value = "null,null,null,null,null,null,null,null,null,Python,null,null"

$classes = $wb.Worksheets.Item("Class Information")
$classes.Cells.Item(2, 7).Value = 1
$classes.Cells.Item(11, 7).Value = 2
$classes.Cells.Item(12, 7).Value = 2
